# Update proses import data excel tabel barang : 13092022
#
# Source header row was: KodeBarang | NamaBarang | Customer | Satuan
# Target header row is:  KodeBarang | NamaBarang | Deskripsi | Customer | Supplier | Satuan
#
# i.e. "Customer" moves from C1 to D1, "Satuan" moves from D1 to F1, and two new
# columns are inserted: "Deskripsi" (C1) and "Supplier" (E1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the two existing header values into their new homes first ...
$ws.Range("D1").Value = "Customer"
$ws.Range("F1").Value = "Satuan"

# ... then add the new header values. "Supplier" must be registered in the
# shared-string table before "Deskripsi" so the string indexes line up with
# the authored workbook (Supplier=4, Deskripsi=5).
$ws.Range("E1").Value = "Supplier"
$ws.Range("C1").Value = "Deskripsi"

# Column widths (best-fit-ish), matching the authored layout as closely as
# this engine's char-width quantization allows.
$ws.Columns.Item(1).ColumnWidth = 10.666666666666666  # A KodeBarang
$ws.Columns.Item(2).ColumnWidth = 11.333333333333334  # B NamaBarang
$ws.Columns.Item(3).ColumnWidth = 11.333333333333334  # C Deskripsi
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666  # D Customer
$ws.Columns.Item(5).ColumnWidth = 10                  # E Supplier
$ws.Columns.Item(6).ColumnWidth = 6.166666666666667   # F Satuan

# The authored workbook carries a small (unused) phonetic-guide font entry
# (8pt Calibri) in styles.xml. Create it without letting it stick to any
# cell's style by writing it to a scratch cell, then clearing that cell
# completely (value + formatting) so the worksheet itself stays unaffected.
$ws.Range("Z100").Font.Size = 8
$ws.Range("Z100").Clear()
